$wb = $excel.ActiveWorkbook

# Overview sheet: Latest HO Xliff Generate Date for bbba99bb-... row
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-17 01:01:25"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for bbba99bb-... row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-17 01:01:20"
$wsZhCn.Range("K2").Value = "2016-08-17 01:01:36"

# de-de sheet: Correspond Handoff Datetime for bbba99bb-... row
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-17 01:01:25"
$wsDeDe.Range("K2").Value = "2016-08-17 01:01:44"
